$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...2017 – on donport4..." -> "...2017 updated 07/24/2018 – on donport4..."
# Original run text: " – on "  -> split into " " / "updated 07/24/2018 " / "– on "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(" " + [char]0x2013 + " on ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insPoint = $d.Range($rng.Start + 1, $rng.Start + 1)
    $insPoint.InsertBefore("updated 07/24/2018 ")
}

# ---------------------------------------------------------------------------
# Change 2: "Using Commandbox " -> "Maintain the installation using Commandbox "
# (bold preserved throughout)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Using Commandbox ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Maintain the installation using Commandbox "
    $rng.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 3: "Here is link to documentation: " -> "Here is link to CFConfig documentation: "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Here is link to documentation: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Here is link to CFConfig documentation: "
    $rng.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 4: Insert "Step one, from scratch, " before "Install cfconfig..." paragraph
# and make "cfconfig" bold
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Install ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insPoint = $d.Range($rng.Start, $rng.Start)
    $insPoint.InsertBefore("Step one, from scratch, ")
}

$rng = $d.Content
$found = $rng.Find.Execute("Step one, from scratch, Install ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $cfStart = $rng.End
    $cfRange = $d.Range($cfStart, $cfStart + 8)
    Write-Host "cfRange text: [" $cfRange.Text "]"
    $cfRange.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 5: "box prompt with" -> "box" (bold) + " prompt with"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("box prompt with", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $boxStart = $rng.Start
    $boxRange = $d.Range($boxStart, $boxStart + 3)
    $boxRange.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 6: "machine,  server" -> "machine, server" (remove double space + proofErr gram markers)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("machine,  server", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "machine, server"
}

# ---------------------------------------------------------------------------
# Change 7: Append new sentences after "...and look like this"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("and look like this", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insPoint = $rng.Duplicate
    $insPoint.Collapse(0)
    $insPoint.InsertAfter(". There are Navsea, and BBSG definitions in there. CF2016 is the Navsea definitions that go with CF 2016.")
}

# ---------------------------------------------------------------------------
# Change 8: "to store the json version of " (merge runs, no visible text change)
# Already represented correctly by plain text - no-op needed since text unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Change 12: Add "SAMPLES:  " before "MAKE SURE THERE ARE NO SPACES..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("MAKE SURE THERE ARE NO SPACES", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insPoint = $d.Range($rng.Start, $rng.Start)
    $insPoint.InsertBefore("SAMPLES:  ")
    $newRange = $d.Range($rng.Start - 10, $rng.Start)
}

# ---------------------------------------------------------------------------
# Change 14/15: BBSG "cfconfig export..." paragraph / "Setting CF admin password"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Setting CF admin password: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Setting password for CF Admin in CF2016: "
}

Write-Host "Edits applied"
